$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency Price (D) and Volume(1h) (E) columns with latest values.
# NumberFormat is set to text ("@") for Price cells first so that numeric-looking
# strings (e.g. "329.16", "0.08890") are kept verbatim instead of being coerced
# into floating point numbers by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.530.37"
$ws.Range("E2").Value = "  +0.54%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.107.40"
$ws.Range("E3").Value = "  +4.78%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "329.16"
$ws.Range("E5").Value = "  +1.26%  "
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5259"
$ws.Range("E7").Value = "  +2.51%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4368"
$ws.Range("E8").Value = "  +2.48%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08890"
$ws.Range("E9").Value = "  +1.97%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "47.35"
$ws.Range("E10").Value = "  +9.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.163"
$ws.Range("E11").Value = "  +2.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.64"
$ws.Range("E12").Value = "  -0.20%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.108.39"
$ws.Range("E13").Value = "  +4.76%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.736"
$ws.Range("E14").Value = "  +2.32%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.766"
$ws.Range("E15").Value = "  +4.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "96.35"
$ws.Range("E16").Value = "  +2.13%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.002"
$ws.Range("E18").Value = "  +1.32%  "
$ws.Range("E19").Value = "  +1.65%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.01"
$ws.Range("E20").Value = "  +0.58%  "
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.320"
$ws.Range("E22").Value = "  +1.87%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.594.37"
$ws.Range("E23").Value = "  +0.59%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.27"
$ws.Range("E24").Value = "  +3.79%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.355"
$ws.Range("E25").Value = "  +4.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.355.41"
$ws.Range("E26").Value = "  +4.78%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.42"
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.609"
$ws.Range("E28").Value = "  +7.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "161.90"
$ws.Range("E29").Value = "  -0.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.78"
$ws.Range("E30").Value = "  +1.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.203"
$ws.Range("E31").Value = "  +5.66%  "
$ws.Range("E32").Value = "  +2.37%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.678"
$ws.Range("E33").Value = "  +22.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.204"
$ws.Range("E34").Value = "  +2.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.927"
$ws.Range("E35").Value = "  +2.61%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.14"
$ws.Range("E36").Value = "  +11.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02583"
$ws.Range("E37").Value = "  +2.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.498"
$ws.Range("E38").Value = "  +0.67%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "12.74"
$ws.Range("E39").Value = "  +2.77%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06687"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2288"
$ws.Range("E41").Value = "  +4.28%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6842"
$ws.Range("E42").Value = "  +2.82%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.263"
$ws.Range("E43").Value = "  +1.89%  "
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.06"
$ws.Range("E45").Value = "  +2.61%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6381"
$ws.Range("E46").Value = "  +3.49%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.209"
$ws.Range("E47").Value = "  +0.79%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.625"
$ws.Range("E48").Value = "  -1.17%  "
$ws.Range("E49").Value = "  -0.71%  "
$ws.Range("E50").Value = "  +7.79%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "82.44"
$ws.Range("E51").Value = "  +2.23%  "